# Updates cryptos list (price / 1h volume / row reorders) per the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '66.914.53'
$ws.Range('E2').Value = '  +0.18%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.822.44'
$ws.Range('E3').Value = '  +3.40%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.18%  '

# Row 5: BNB
$ws.Range('D5').Value = '''412.78'
$ws.Range('E5').Value = '  -1.57%  '

# Row 6: Solana
$ws.Range('D6').Value = '''132.64'
$ws.Range('E6').Value = '  +1.79%  '

# Row 7: LidoStakedEther
$ws.Range('D7').Value = '3.811.42'
$ws.Range('E7').Value = '  +3.30%  '

# Row 8: XRP
$ws.Range('D8').Value = '''0.616'
$ws.Range('E8').Value = '  -4.09%  '

# Row 9: USDC
$ws.Range('E9').Value = '  +0.07%  '

# Row 10: Cardano
$ws.Range('D10').Value = '''0.741'
$ws.Range('E10').Value = '  -2.74%  '

# Row 11: Dogecoin
$ws.Range('D11').Value = '''0.171'
$ws.Range('E11').Value = '  -5.62%  '

# Row 12: ShibaInu
$ws.Range('D12').Value = '''0.0000379'
$ws.Range('E12').Value = '  -3.36%  '

# Row 13: Avalanche
$ws.Range('D13').Value = '''41.16'
$ws.Range('E13').Value = '  -4.85%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('D14').Value = '4.443.29'
$ws.Range('E14').Value = '  +3.82%  '

# Row 15: Polkadot
$ws.Range('D15').Value = '''10.02'
$ws.Range('E15').Value = '  -6.33%  '

# Row 16: Uniswap
$ws.Range('D16').Value = '''14.87'
$ws.Range('E16').Value = '  +13.29%  '

# Row 17: TRON
$ws.Range('E17').Value = '  -1.14%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '3.822.20'
$ws.Range('E18').Value = '  +0.30%  '

# Row 19: Chainlink
$ws.Range('D19').Value = '''19.53'
$ws.Range('E19').Value = '  -4.81%  '

# Row 20: WrappedBTC
$ws.Range('D20').Value = '67.236.06'
$ws.Range('E20').Value = '  +0.71%  '

# Row 21: Polygon
$ws.Range('E21').Value = '  -2.95%  '

# Row 22: BitcoinCash
$ws.Range('D22').Value = '''415.34'
$ws.Range('E22').Value = '  -6.15%  '

# Row 23: InternetComputer(DFINITY)
$ws.Range('D23').Value = '''14.91'
$ws.Range('E23').Value = '  -8.19%  '

# Row 24: Litecoin
$ws.Range('D24').Value = '''86.15'
$ws.Range('E24').Value = '  -3.94%  '

# Row 25: ImmutableX
$ws.Range('D25').Value = '''3.09'
$ws.Range('E25').Value = '  -1.54%  '

# Row 26: EthereumClassic
$ws.Range('D26').Value = '''36.54'
$ws.Range('E26').Value = '  -2.23%  '

# Row 27: LEO
$ws.Range('E27').Value = '  +14.09%  '

# Row 28: PancakeSwap
$ws.Range('D28').Value = '''3.15'
$ws.Range('E28').Value = '  -4.61%  '

# Row 29: Filecoin
$ws.Range('D29').Value = '''9.53'
$ws.Range('E29').Value = '  -7.08%  '

# Row 30: Bittensor
$ws.Range('D30').Value = '''691.13'
$ws.Range('E30').Value = '  +5.95%  '

# Row 31: Hedera
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.122'
$ws.Range('E31').Value = '  -2.50%  '

# Row 32: Cosmos
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '''12.49'
$ws.Range('E32').Value = '  -2.22%  '

# Row 33: Toncoin
$ws.Range('E33').Value = '  +1.07%  '

# Row 34: RenderToken
$ws.Range('D34').Value = '''7.26'
$ws.Range('E34').Value = '  -0.60%  '

# Row 35: Kaspa
$ws.Range('E35').Value = '  -7.48%  '

# Row 36: InjectiveProtocol
$ws.Range('D36').Value = '''39.12'
$ws.Range('E36').Value = '  -6.46%  '

# Row 37: Dai
$ws.Range('E37').Value = '  -0.08%  '

# Row 38: OKB
$ws.Range('D38').Value = '''55.50'
$ws.Range('E38').Value = '  -2.89%  '

# Row 39: PEPE
$ws.Range('D39').Value = '0.0₃0778'
$ws.Range('E39').Value = '  +6.14%  '

# Row 40: VeChain
$ws.Range('D40').Value = '''0.0462'
$ws.Range('E40').Value = '  -6.38%  '

# Row 41: ThetaToken
$ws.Range('D41').Value = '''3.06'
$ws.Range('E41').Value = '  -0.67%  '

# Row 42: FirstDigitalUSD
$ws.Range('D42').Value = '''0.999'
$ws.Range('E42').Value = '  +0.12%  '

# Row 43: EnergySwap
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '''27.71'
$ws.Range('E43').Value = '  -4.73%  '

# Row 44: Stellar
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '''0.137'
$ws.Range('E44').Value = '  -8.62%  '

# Row 45: Monero
$ws.Range('D45').Value = '''148.35'
$ws.Range('E45').Value = '  -0.25%  '

# Row 46: ApeXProtocol
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '''3.18'
$ws.Range('E46').Value = '  +19.43%  '

# Row 47: LidoDAOToken
$ws.Range('B47').Value = 'LidoDAOToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D47').Value = '''3.33'
$ws.Range('E47').Value = '  -2.19%  '

# Row 48: NEARProtocol
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''4.42'
$ws.Range('E48').Value = '  +1.64%  '

# Row 49: ARBITRUM
$ws.Range('D49').Value = '''2.10'
$ws.Range('E49').Value = '  -0.47%  '

# Row 50: Stacks
$ws.Range('D50').Value = '''2.86'
$ws.Range('E50').Value = '  -1.00%  '

# Row 51: WEMIXToken
$ws.Range('D51').Value = '''2.60'
$ws.Range('E51').Value = '  -1.34%  '
